$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($r, $c, $v) {
    $ws.Cells.Item($r, $c).Value = $v
}

# Header renames (column headers -> snake_case machine names)
Set-Cell 1 1 "mx_state"
Set-Cell 1 2 "mx_municipality"
Set-Cell 1 3 "n_matriculas"
Set-Cell 1 4 "pct_matriculas"

# Title-case the Spanish connector words (de/del/el/la/los/las/y) in place names
Set-Cell 5 2 "Pabellón De Arteaga"
Set-Cell 6 2 "Rincón De Romos"
Set-Cell 7 2 "San José De Gracia"
Set-Cell 24 2 "Amatenango De La Frontera"
Set-Cell 32 2 "Comitán De Domínguez"
Set-Cell 46 2 "Mazapa De Madero"
Set-Cell 68 2 "Guadalupe Y Calvo"
Set-Cell 70 2 "Hidalgo Del Parral"
Set-Cell 98 2 "Villa De Álvarez"
Set-Cell 100 1 "Ciudad De México"
Set-Cell 128 2 "San Juan Del Río"
Set-Cell 136 1 "Estado De México"
Set-Cell 136 2 "Acambay De Ruíz Castañeda"
Set-Cell 139 2 "Almoloya De Juárez"
Set-Cell 143 2 "Atizapán De Zaragoza"
Set-Cell 148 2 "Chapa De Mota"
Set-Cell 150 2 "Coacalco De Berriozábal"
Set-Cell 155 2 "Ecatepec De Morelos"
Set-Cell 159 2 "Ixtapan De La Sal"
Set-Cell 168 2 "Naucalpan De Juárez"
Set-Cell 177 2 "San Felipe Del Progreso"
Set-Cell 187 2 "Tenango Del Valle"
Set-Cell 197 2 "Tlalnepantla De Baz"
Set-Cell 201 2 "Valle De Bravo"
Set-Cell 202 2 "Valle De Chalco Solidaridad"
Set-Cell 203 2 "Villa De Allende"
Set-Cell 204 2 "Villa Del Carbón"
Set-Cell 215 2 "San Miguel De Allende"
Set-Cell 216 2 "Apaseo El Alto"
Set-Cell 217 2 "Apaseo El Grande"
Set-Cell 224 2 "Dolores Hidalgo Cuna De La Independencia Nacional"
Set-Cell 227 2 "Jaral Del Progreso"
Set-Cell 234 2 "Purísima Del Rincón"
Set-Cell 239 2 "San Francisco Del Rincón"
Set-Cell 241 2 "San Luis De La Paz"
Set-Cell 242 2 "Silao De La Victoria"
Set-Cell 246 2 "Valle De Santiago"
Set-Cell 252 2 "Acapulco De Juárez"
Set-Cell 254 2 "Ajuchitlán Del Progreso"
Set-Cell 255 2 "Alcozauca De Guerrero"
Set-Cell 259 2 "Atenango Del Río"
Set-Cell 261 2 "Atoyac De Álvarez"
Set-Cell 262 2 "Ayutla De Los Libres"
Set-Cell 265 2 "Buenavista De Cuéllar"
Set-Cell 266 2 "Chilapa De Álvarez"
Set-Cell 267 2 "Chilpancingo De Los Bravo"
Set-Cell 268 2 "Coahuayutla De José María Izazaga"
Set-Cell 273 2 "Coyuca De Benítez"
Set-Cell 274 2 "Coyuca De Catalán"
Set-Cell 278 2 "Cuetzala Del Progreso"
Set-Cell 279 2 "Cutzamala De Pinzón"
Set-Cell 284 2 "Huitzuco De Los Figueroa"
Set-Cell 285 2 "Iguala De La Independencia"
Set-Cell 287 2 "Ixcateopan De Cuauhtémoc"
Set-Cell 288 2 "Zihuatanejo De Azueta"
Set-Cell 290 2 "La Unión De Isidoro Montes De Oca"
Set-Cell 293 2 "Mártir De Cuilapan"
Set-Cell 306 2 "Taxco De Alarcón"
Set-Cell 308 2 "Técpan De Galeana"
Set-Cell 310 2 "Tepecoacuilco De Trujano"
Set-Cell 311 2 "Tixtla De Guerrero"
Set-Cell 315 2 "Tlalixtaquilla De Maldonado"
Set-Cell 316 2 "Tlapa De Comonfort"
Set-Cell 328 2 "Atotonilco El Grande"
Set-Cell 332 2 "Cuautepec De Hinojosa"
Set-Cell 336 2 "Huasca De Ocampo"
Set-Cell 337 2 "Huejutla De Reyes"
Set-Cell 340 2 "Jacala De Ledezma"
Set-Cell 345 2 "Mineral Del Monte"
Set-Cell 346 2 "Mixquiahuala De Juárez"
Set-Cell 347 2 "Molango De Escamilla"
Set-Cell 348 2 "Nopala De Villagrán"
Set-Cell 349 2 "Omitlán De Juárez"
Set-Cell 350 2 "Pachuca De Soto"
Set-Cell 352 2 "Progreso De Obregón"
Set-Cell 356 2 "Santiago De Anaya"
Set-Cell 357 2 "Santiago Tulantepec De Lugo Guerrero"
Set-Cell 361 2 "Tenango De Doria"
Set-Cell 363 2 "Tezontepec De Aldama"
Set-Cell 368 2 "Tula De Allende"
Set-Cell 369 2 "Tulancingo De Bravo"
Set-Cell 371 2 "Zapotlán De Juárez"
Set-Cell 375 2 "Acatlán De Juárez"
Set-Cell 376 2 "Ahualulco De Mercado"
Set-Cell 381 2 "Atemajac De Brizuela"
Set-Cell 383 2 "Atotonilco El Alto"
Set-Cell 385 2 "Autlán De Navarro"
Set-Cell 395 2 "Concepción De Buenos Aires"
Set-Cell 396 2 "Cuautitlán De García Barragán"
Set-Cell 409 2 "Huejuquilla El Alto"
Set-Cell 410 2 "Ixtlahuacán De Los Membrillos"
Set-Cell 411 2 "Ixtlahuacán Del Río"
Set-Cell 415 2 "Jilotlán De Los Dolores"
Set-Cell 421 2 "La Manzanilla De La Paz"
Set-Cell 422 2 "Lagos De Moreno"
Set-Cell 429 2 "Ojuelos De Jalisco"
Set-Cell 434 2 "San Diego De Alejandría"
Set-Cell 436 2 "San Juan De Los Lagos"
Set-Cell 439 2 "San Martín De Bolaños"
Set-Cell 441 2 "San Miguel El Alto"
Set-Cell 442 2 "San Sebastián Del Oeste"
Set-Cell 443 2 "Santa María De Los Ángeles"
Set-Cell 446 2 "Talpa De Allende"
Set-Cell 447 2 "Tamazula De Gordiano"
Set-Cell 451 2 "Teocuitatlán De Corona"
Set-Cell 452 2 "Tepatitlán De Morelos"
Set-Cell 455 2 "Tizapán El Alto"
Set-Cell 456 2 "Tlajomulco De Zúñiga"
Set-Cell 467 2 "Unión De San Antonio"
Set-Cell 468 2 "Unión De Tula"
Set-Cell 472 2 "Yahualica De González Gallo"
Set-Cell 473 2 "Zacoalco De Torres"
Set-Cell 476 2 "Zapotitlán De Vadillo"
Set-Cell 477 2 "Zapotlán Del Rey"
Set-Cell 478 2 "Zapotlán El Grande"
Set-Cell 502 2 "Coalcomán De Vázquez Pallares"
Set-Cell 504 2 "Cojumatlán De Régules"
Set-Cell 566 2 "Tiquicheo De Nicolás Romero"
Set-Cell 596 2 "Jonacatepec De Leandro Valle"
Set-Cell 599 2 "Puente De Ixtla"
Set-Cell 602 2 "Tetela Del Volcán"
Set-Cell 603 2 "Tlaltizapán De Zapata"
Set-Cell 611 2 "Zacualpan De Amilpas"
Set-Cell 615 2 "Amatlán De Cañas"
Set-Cell 616 2 "Bahía De Banderas"
Set-Cell 620 2 "Ixtlán Del Río"
Set-Cell 627 2 "Santa María Del Oro"
Set-Cell 638 2 "San Nicolás De Los Garza"
Set-Cell 641 2 "Acatlán De Pérez Figueroa"
Set-Cell 645 2 "Ayoquezco De Aldama"
Set-Cell 648 2 "Chalcatongo De Hidalgo"
Set-Cell 651 2 "Coicoyán De Las Flores"
Set-Cell 652 2 "Constancia Del Rosario"
Set-Cell 655 2 "Cuilápam De Guerrero"
Set-Cell 656 2 "Guadalupe De Ramírez"
Set-Cell 657 2 "Heroica Ciudad De Ejutla De Crespo"
Set-Cell 658 2 "Heroica Ciudad De Huajuapan De León"
Set-Cell 659 2 "Heroica Ciudad De Tlaxiaco"
Set-Cell 660 2 "Huautla De Jiménez"
Set-Cell 662 2 "Ixtlán De Juárez"
Set-Cell 663 2 "Heroica Ciudad De Juchitán De Zaragoza"
Set-Cell 672 2 "Mariscala De Juárez"
Set-Cell 675 2 "Miahuatlán De Porfirio Díaz"
Set-Cell 676 2 "Mixistlán De La Reforma"
Set-Cell 678 2 "Oaxaca De Juárez"
Set-Cell 679 2 "Ocotlán De Morelos"
Set-Cell 680 2 "Pinotepa De Don Luis"
Set-Cell 681 2 "Putla Villa De Guerrero"
Set-Cell 693 2 "San Baltazar Yatzachi El Bajo"
Set-Cell 698 2 "San Felipe Jalapa De Díaz"
Set-Cell 709 2 "San José Del Peñasco"
Set-Cell 714 2 "San Juan Bautista Lo De Soto"
Set-Cell 721 2 "San Juan Del Estado"
Set-Cell 722 2 "San Juan Del Río"
Set-Cell 745 2 "San Mateo Del Mar"
Set-Cell 756 2 "San Miguel Del Puerto"
Set-Cell 758 2 "San Miguel El Grande"
Set-Cell 771 2 "San Pablo Villa De Mitla"
Set-Cell 773 2 "San Pedro El Alto"
Set-Cell 781 2 "San Pedro Y San Pablo Teposcolula"
Set-Cell 798 2 "Santa Cruz Tacache De Mina"
Set-Cell 800 2 "Santa Inés De Zaragoza"
Set-Cell 801 2 "Santa Inés Del Monte"
Set-Cell 813 2 "Santa María Jalapa Del Marqués"
Set-Cell 853 2 "Santo Domingo De Morelos"
Set-Cell 866 2 "Sitio De Xitlapehua"
Set-Cell 868 2 "Tamazulápam Del Espíritu Santo"
Set-Cell 870 2 "Tataltepec De Valdés"
Set-Cell 871 2 "Teococuilco De Marcos Pérez"
Set-Cell 872 2 "Teotitlán De Flores Magón"
Set-Cell 873 2 "Teotitlán Del Valle"
Set-Cell 875 2 "Tepelmeme Villa De Morelos"
Set-Cell 876 2 "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"
Set-Cell 877 2 "Tlacolula De Matamoros"
Set-Cell 878 2 "Tlalixtac De Cabrera"
Set-Cell 881 2 "Villa De Etla"
Set-Cell 882 2 "Villa De Tututepec"
Set-Cell 883 2 "Villa De Zaachila"
Set-Cell 886 2 "Villa Sola De Vega"
Set-Cell 888 2 "Zimatlán De Álvarez"
Set-Cell 903 2 "Chalchicomula De Sesma"
Set-Cell 911 2 "Chila De La Sal"
Set-Cell 918 2 "Cuayuca De Andrade"
Set-Cell 919 2 "Cuetzalan Del Progreso"
Set-Cell 927 2 "Huehuetlán El Chico"
Set-Cell 928 2 "Huehuetlán El Grande"
Set-Cell 934 2 "Izúcar De Matamoros"
Set-Cell 938 2 "Los Reyes De Juárez"
Set-Cell 945 2 "Palmar De Bravo"
Set-Cell 960 2 "San Salvador El Seco"
Set-Cell 964 2 "Tecali De Herrera"
Set-Cell 971 2 "Tepanco De López"
Set-Cell 974 2 "Tepexi De Rodríguez"
Set-Cell 976 2 "Tetela De Ocampo"
Set-Cell 980 2 "Tlacotepec De Benito Juárez"
Set-Cell 988 2 "Totoltepec De Guerrero"
Set-Cell 993 2 "Xayacatlán De Bravo"
Set-Cell 996 2 "Xochitlán De Vicente Suárez"
Set-Cell 1004 2 "Amealco De Bonfil"
Set-Cell 1006 2 "Cadereyta De Montes"
Set-Cell 1011 2 "Pinal De Amoles"
Set-Cell 1014 2 "San Juan Del Río"
Set-Cell 1022 2 "Armadillo De Los Infante"
Set-Cell 1029 2 "Mexquitic De Carmona"
Set-Cell 1032 2 "San Ciro De Acosta"
Set-Cell 1035 2 "Santa María Del Río"
Set-Cell 1037 2 "Soledad De Graciano Sánchez"
Set-Cell 1040 2 "Villa De Ramos"
Set-Cell 1092 2 "Soto La Marina"
Set-Cell 1097 2 "Acuamanala De Miguel Hidalgo"
Set-Cell 1103 2 "Ixtacuixtla De Mariano Matamoros"
Set-Cell 1108 2 "San Pablo Del Monte"
Set-Cell 1109 2 "Sanctórum De Lázaro Cárdenas"
Set-Cell 1112 2 "Tepetitla De Lardizábal"
Set-Cell 1127 2 "Amatlán De Los Reyes"
Set-Cell 1138 2 "Cazones De Herrera"
Set-Cell 1144 2 "Cosamaloapan De Carpio"
Set-Cell 1155 2 "Hueyapan De Ocampo"
Set-Cell 1156 2 "Huiloapan De Cuauhtémoc"
Set-Cell 1159 2 "Ixhuatlán De Madero"
Set-Cell 1160 2 "Ixhuatlán Del Café"
Set-Cell 1161 2 "Ixhuatlán Del Sureste"
Set-Cell 1175 2 "Martínez De La Torre"
Set-Cell 1176 2 "Medellín De Bravo"
Set-Cell 1187 2 "Ozuluama De Mascareñas"
Set-Cell 1190 2 "Paso De Ovejas"
Set-Cell 1192 2 "Poza Rica De Hidalgo"
Set-Cell 1198 2 "Sayula De Alemán"
Set-Cell 1201 2 "Soledad De Doblado"
Set-Cell 1222 2 "Vega De Alatorre"
Set-Cell 1253 2 "Cañitas De Felipe Pescador"
Set-Cell 1264 2 "Mezquital Del Oro"
Set-Cell 1267 2 "Moyahua De Estrada"
Set-Cell 1268 2 "Nochistlán De Mejía"
Set-Cell 1269 2 "Noria De Ángeles"
Set-Cell 1277 2 "Teúl De González Ortega"
Set-Cell 1278 2 "Tlaltenango De Sánchez Román"
Set-Cell 1280 2 "Villa De Cos"

# Remove trailing footnote/metadata rows (1286:1290); the sheet dimension
# auto-shrinks to the new used range (A1:D1284) once they are cleared.
$ws.Range("A1286:D1290").ClearContents()
